$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price/Volume columns first so values such as
# "29.40", "4.10" or "1.622.65" stay literal strings instead of being
# auto-coerced into numbers by Excel (matching the source inlineStr cells).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '29.866.84'
$ws.Cells.Item(2, 5).Value = '  +1.18%  '
$ws.Cells.Item(3, 4).Value = '1.622.65'
$ws.Cells.Item(3, 5).Value = '  +1.21%  '
$ws.Cells.Item(4, 4).Value = '0.994'
$ws.Cells.Item(4, 5).Value = '  -0.47%  '
$ws.Cells.Item(5, 4).Value = '213.35'
$ws.Cells.Item(5, 5).Value = '  +0.46%  '
$ws.Cells.Item(6, 4).Value = '0.518'
$ws.Cells.Item(6, 5).Value = '  -0.78%  '
$ws.Cells.Item(7, 4).Value = '0.994'
$ws.Cells.Item(7, 5).Value = '  -0.41%  '
$ws.Cells.Item(8, 4).Value = '29.40'
$ws.Cells.Item(8, 5).Value = '  +9.64%  '
$ws.Cells.Item(9, 5).Value = '  +3.18%  '
$ws.Cells.Item(10, 5).Value = '  +0.89%  '
$ws.Cells.Item(11, 5).Value = '  +0.18%  '
$ws.Cells.Item(12, 4).Value = '1.855.69'
$ws.Cells.Item(12, 5).Value = '  +1.24%  '
$ws.Cells.Item(13, 4).Value = '1.618.65'
$ws.Cells.Item(13, 5).Value = '  +1.16%  '
$ws.Cells.Item(14, 4).Value = '0.567'
$ws.Cells.Item(14, 5).Value = '  +5.81%  '
$ws.Cells.Item(15, 5).Value = '  +5.70%  '
$ws.Cells.Item(16, 4).Value = '29.921.46'
$ws.Cells.Item(16, 5).Value = '  +1.39%  '
$ws.Cells.Item(17, 4).Value = '8.81'
$ws.Cells.Item(17, 5).Value = '  +15.90%  '
$ws.Cells.Item(18, 4).Value = '64.29'
$ws.Cells.Item(18, 5).Value = '  +1.30%  '
$ws.Cells.Item(19, 4).Value = '242.12'
$ws.Cells.Item(19, 5).Value = '  +0.72%  '
$ws.Cells.Item(20, 5).Value = '  +2.28%  '
$ws.Cells.Item(21, 5).Value = '  -0.32%  '
$ws.Cells.Item(22, 4).Value = '4.10'
$ws.Cells.Item(22, 5).Value = '  +2.98%  '
$ws.Cells.Item(23, 5).Value = '  +4.16%  '
$ws.Cells.Item(24, 4).Value = '2.14'
$ws.Cells.Item(24, 5).Value = '  +2.73%  '
$ws.Cells.Item(25, 4).Value = '156.68'
$ws.Cells.Item(25, 5).Value = '  +1.42%  '
$ws.Cells.Item(26, 4).Value = '15.60'
$ws.Cells.Item(26, 5).Value = '  +2.21%  '
$ws.Cells.Item(27, 4).Value = '0.111'
$ws.Cells.Item(27, 5).Value = '  +1.43%  '
$ws.Cells.Item(29, 5).Value = '  -0.42%  '
$ws.Cells.Item(30, 4).Value = '0.0488'
$ws.Cells.Item(30, 5).Value = '  +3.31%  '
$ws.Cells.Item(31, 5).Value = '  +5.44%  '
$ws.Cells.Item(32, 4).Value = '3.33'
$ws.Cells.Item(32, 5).Value = '  +3.31%  '
$ws.Cells.Item(33, 5).Value = '  +4.29%  '
$ws.Cells.Item(34, 4).Value = '1.424.05'
$ws.Cells.Item(34, 5).Value = '  +1.06%  '
$ws.Cells.Item(35, 4).Value = '1.64'
$ws.Cells.Item(35, 5).Value = '  +6.94%  '
$ws.Cells.Item(36, 5).Value = '  -0.88%  '
$ws.Cells.Item(37, 4).Value = '2.87'
$ws.Cells.Item(37, 5).Value = '  +1.90%  '
$ws.Cells.Item(38, 4).Value = '2.29'
$ws.Cells.Item(38, 5).Value = '  -0.64%  '
$ws.Cells.Item(39, 5).Value = '  +2.87%  '
$ws.Cells.Item(40, 4).Value = '0.556'
$ws.Cells.Item(40, 5).Value = '  +3.29%  '
$ws.Cells.Item(41, 5).Value = '  +3.49%  '
$ws.Cells.Item(42, 5).Value = '  -0.15%  '
$ws.Cells.Item(43, 4).Value = '0.827'
$ws.Cells.Item(43, 5).Value = '  +3.75%  '
$ws.Cells.Item(44, 4).Value = '53.74'
$ws.Cells.Item(44, 5).Value = '  +1.44%  '
$ws.Cells.Item(45, 4).Value = '69.15'
$ws.Cells.Item(45, 5).Value = '  +5.12%  '
$ws.Cells.Item(46, 5).Value = '  +18.82%  '
$ws.Cells.Item(47, 4).Value = '0.993'
$ws.Cells.Item(47, 5).Value = '  -0.48%  '
$ws.Cells.Item(48, 5).Value = '  +2.78%  '
$ws.Cells.Item(49, 4).Value = '1.764.67'
$ws.Cells.Item(49, 5).Value = '  +1.30%  '
$ws.Cells.Item(50, 4).Value = '88.13'
$ws.Cells.Item(50, 5).Value = '  +1.95%  '
$ws.Cells.Item(51, 5).Value = '  +7.14%  '

# Reset the style back to Normal so we do not leave a stray cell style
# behind -- only the cell values should differ from the source workbook.
$ws.Range("D2:E51").Style = "Normal"
